$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 103 (pushes the old row 103 "Red (Taylor's Version)" / 44512 down to row 104).
# The inserted row inherits styling from the surrounding rows (B: text style, C: date style).
$ws.Rows.Item(103).Insert()

# Add a further new row at the end (row 105) for "This Love (Taylor's Version)".
# Setting the value here first means it gets the earlier shared-string index,
# matching the order the strings were added to xl/sharedStrings.xml.
$ws.Range("B105").Value = "This Love (Taylor's Version)"
$ws.Range("C105").Value = 44687

# Fill in the new row 103 with "The Joker And The Queen".
$ws.Range("B103").Value = "The Joker And The Queen"
$ws.Range("C103").Value = 44498

# Restore the selection to match the edited view state.
$ws.Range("B103").Select() | Out-Null
